$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.8     # Current Capital
$summary.Range("B4").Value = -0.2       # Total P&L $
$summary.Range("B5").Value = -0.57      # Total P&L %
$summary.Range("B6").Value = 7          # Total Trades
$summary.Range("B7").Value = 3          # Winning Trades
$summary.Range("B9").Value = 42.86      # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.8
$status.Range("D5").Value = 7
$status.Range("E5").Value = -0.2
$status.Range("F5").Value = -0.2
$status.Range("G5").Value = 42.86

# ---------------------------------------------------------------------------
# Append new trade row (#7) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
# Columns C.. onward (skip B - the "Date" column is handled separately below,
# since assigning a date-shaped literal like "2026-02-17" to .Value gets
# auto-coerced into a date serial number by the engine, same as real Excel).
$newRowRest = @{
    1  = 7                                         # A: Trade #
    3  = "20:02:28"                                # C: Time
    4  = "MarketMaking"                            # D: Strategy
    5  = "UP"                                      # E: Side
    6  = 0.57                                      # F: Entry Price
    7  = 0.613861                                  # G: Exit Price
    8  = "CLOSED"                                  # H: Status
    9  = 7.695                                     # I: P&L %
    10 = 0.04                                      # J: P&L $
    11 = 99.8                                      # K: Capital After
    12 = 0                                         # L: Entry Slippage (bps)
    13 = 0                                         # M: Exit Slippage (bps)
    14 = 0.6                                       # N: Confidence
    15 = "Normal spread capture: 19600 bps"        # O: Entry Reason
    16 = "early_exit"                              # P: Exit Reason
    17 = 0.13                                      # Q: Duration (min)
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 8

    # "Date" column: copy from an existing row's date cell (same literal
    # text "2026-02-17") so it is written as plain text, not re-parsed into
    # a date serial number.
    $ws.Range("B2").Copy($ws.Cells.Item($row, 2))

    foreach ($col in $newRowRest.Keys) {
        $ws.Cells.Item($row, $col).Value = $newRowRest[$col]
    }
}
